$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Diagonal cells (origin == destination) representing intra-municipality trips.
# Values are being increased: kids moved from the "non_work" matrix into the "working" matrix.
$updates = @(
    @{Cell="B2"; Value=5033}
    @{Cell="C3"; Value=42870}
    @{Cell="D4"; Value=9109}
    @{Cell="E5"; Value=3810}
    @{Cell="F6"; Value=9758}
    @{Cell="G7"; Value=43969}
    @{Cell="H8"; Value=17408}
    @{Cell="I9"; Value=15910}
    @{Cell="J10"; Value=11215}
    @{Cell="K11"; Value=6167}
    @{Cell="L12"; Value=19468}
    @{Cell="M13"; Value=18114}
    @{Cell="N14"; Value=123413}
    @{Cell="O15"; Value=11473}
    @{Cell="P16"; Value=16711}
    @{Cell="Q17"; Value=19046}
    @{Cell="R18"; Value=13125}
    @{Cell="S19"; Value=14181}
    @{Cell="T20"; Value=9461}
    @{Cell="U21"; Value=160999}
    @{Cell="V22"; Value=24970}
    @{Cell="W23"; Value=12093}
    @{Cell="X24"; Value=19875}
    @{Cell="Y25"; Value=17581}
    @{Cell="Z26"; Value=13330}
    @{Cell="AA27"; Value=11306}
    @{Cell="AB28"; Value=7898}
    @{Cell="AC29"; Value=7291}
    @{Cell="AD30"; Value=15755}
    @{Cell="AE31"; Value=8269}
    @{Cell="AF32"; Value=57838}
    @{Cell="AG33"; Value=25406}
    @{Cell="AH34"; Value=7413}
    @{Cell="AI35"; Value=11973}
    @{Cell="AJ36"; Value=8548}
    @{Cell="AK37"; Value=11487}
    @{Cell="AL38"; Value=60476}
    @{Cell="AM39"; Value=5966}
    @{Cell="AN40"; Value=23004}
    @{Cell="AO41"; Value=17271}
    @{Cell="AP42"; Value=9617}
    @{Cell="AQ43"; Value=13177}
    @{Cell="AR44"; Value=30959}
    @{Cell="AS45"; Value=61219}
    @{Cell="AT46"; Value=6360}
    @{Cell="AU47"; Value=47886}
    @{Cell="AV48"; Value=11563}
    @{Cell="AW49"; Value=124178}
    @{Cell="AX50"; Value=3043}
    @{Cell="AY51"; Value=18233}
    @{Cell="AZ52"; Value=11309}
    @{Cell="BA53"; Value=7745}
    @{Cell="BB54"; Value=-15712}
    @{Cell="BC55"; Value=174115}
    @{Cell="BD56"; Value=25378}
    @{Cell="BE57"; Value=11026}
    @{Cell="BF58"; Value=30656}
    @{Cell="BG59"; Value=16027}
    @{Cell="BH60"; Value=3131}
    @{Cell="BI61"; Value=34479}
    @{Cell="BJ62"; Value=11901}
    @{Cell="BK63"; Value=4715}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
